# Progress.xlsx update: log a new entry (row 4) for the "Table1" tracker
# and move the active-cell selection down to the next empty description cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: Date / Time start / Time end for row 4 of Table1
$ws.Range("A4").Value = 45918                     # 2025-09-18
$ws.Range("B4").Value = 0.54166666666666663       # 13:00
$ws.Range("C4").Value = 0.70833333333333337       # 17:00

# Description of the work done in this session (new shared string)
$ws.Range("E4").Value = "Finished Chapter 5. Learned more about debugging. This is something I was not as familiar with as I thought so it was a good thing I learned some proper practices, and Python specific ones too. I alsmost completely finished Chapter 6, except for the last programming exercise. I will finish this next time, this seems like a first proper assignment. This chapter is about learning to use lists better."

# Row grows to fit the wrapped description text (matches row 3's height)
$ws.Rows.Item(4).RowHeight = 72

# Move the selection to the next entry's description cell
$ws.Range("E5").Select()
